$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the three rows that no longer belong on the page-metadata sheet
#    (MechaCar Statistics, Bike Sharing, Employee Database - rows 11-13).
#    Their hyperlink objects have to be removed explicitly first since
#    clearing a range does not drop the workbook's Hyperlinks entries.
# ---------------------------------------------------------------------------
$addressesToRemove = @(
    "https://github.com/cdpeters/MechaCar-statistical-analysis-R",
    "https://github.com/cdpeters/bike-sharing-tableau",
    "https://github.com/cdpeters/employee-database-postgresql",
    "https://public.tableau.com/views/NYC_CitiBike_Challenge_16506220556720/August2019NYCCitibikeStudy?:language=en-US&:display_count=n&:origin=viz_share_link"
)

foreach ($addr in $addressesToRemove) {
    $found = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Address -eq $addr) {
            $found = $h
            break
        }
    }
    if ($found -ne $null) {
        $found.Delete()
    }
}

$ws.Rows("11:13").Clear()

# ---------------------------------------------------------------------------
# 2. Add the new "website_name" column (J) so the project summary row can
#    show a link label whenever a website exists for that project.
# ---------------------------------------------------------------------------
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "website_name"

# Populate the website name for every row that already has a website (col I).
$ws.Range("J16").Value = $ws.Range("D16").Value2
$ws.Range("J17").Value = $ws.Range("D17").Value2

$ws.Columns("J").AutoFit()

# ---------------------------------------------------------------------------
# 3. Restore the active selection to where the author left off editing.
# ---------------------------------------------------------------------------
$ws.Range("A13:XFD13").Select()
